$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "304.34"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "5.80%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "34.90"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "12.28%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.165"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "4.96%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07762"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "6.09%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.290"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.24%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.040"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.90%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.007"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "7.72%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9289"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.84%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1017"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "11.23%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1834"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "8.63%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08624"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "5.35%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03463"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "11.10%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09859"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.71%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001492"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.31%"
$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04620"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2.36%"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005828"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2.28%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.507"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.30%"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.108"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.96%"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3421"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.74%"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1331"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.55%"
$ws.Range("B22").Value = "MCDex"
$ws.Range("C22").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.602"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "9.46%"
$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.2299"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "9.43%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001228"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.38%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004423"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "6.36%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001303"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.24%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003420"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.74%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01768"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "12.44%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04726"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "6.43%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007611"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.87%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "5.77%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007099"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-25.69%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.34%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009201"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "10.48%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00005905"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-3.24%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000752"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.25%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "12.55%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002705"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "35.15%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002106"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.25%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002005"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.25%"
